# Applies the edits described by the diff between the original
# "planilhas_PDA.xlsx" and the updated version uploaded by the author.
#
# Summary of changes:
#  Sheet "PRODUTOS - AULA 1":
#    - B3: 100 -> 500
#    - C3: 100 -> 600
#    - D3: 100 -> 1000
#    - H3: 5 -> 7
#    - L3: 50 -> 60
#    - P3: 50 -> 100
#    (downstream formulas E3, I3, M3, Q3 recalc automatically)
#    - column D becomes wider (content no longer fits the shared 10.33 width)
#    - cursor/selection moved to M3
#  Sheet "SALDO - AULA 2":
#    - D14 formula changed from =C15 to =B15 (and the formulas in
#      E14:G14, which referenced D15/E15/F15, shift accordingly to
#      C15/D15/E15 so the whole row keeps the same "two rows up, one
#      column left" pattern) -- downstream rows 15/16 recalc automatically
#    - cursor/selection moved to G19

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: PRODUTOS - AULA 1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PRODUTOS - AULA 1")

$ws1.Range("B3").Value = 500
$ws1.Range("C3").Value = 600
$ws1.Range("D3").Value = 1000
$ws1.Range("H3").Value = 7
$ws1.Range("L3").Value = 60
$ws1.Range("P3").Value = 100

# Column D no longer fits inside the shared "B:D" 10.33-wide column now
# that it holds a 4-digit number; widen it like Excel's own best-fit did.
$ws1.Columns.Item(4).ColumnWidth = 11

# ---------------------------------------------------------------------
# Sheet 2: SALDO - AULA 2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SALDO - AULA 2")

$ws2.Range("D14").Formula = "=B15"
$ws2.Range("E14:G14").Formula = "=C15"

# ---------------------------------------------------------------------
# Restore selections to match where the author ended up, making sure
# the last-activated sheet stays "SALDO - AULA 2" (it is the tab that
# was active when the file was saved).
# ---------------------------------------------------------------------
$ws1.Range("M3").Select()
$ws2.Range("G19").Select()
$ws2.Activate()
